# Applies the mapping-table fill-in described in the commit diff to the
# "DetectedIssue" worksheet of Profile_DetectedIssue.xlsx.
#
# Column layout (row 2 headers):
#   A = Path
#   B = HL7 v2 Mapping as per Argonaut/www.hl7.org
#   C = HL7 Field (Numeric)
#   D = HL7 Field Name
#   E = Complexity/Operation
#   F = Comments

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F - Comments (filled in for a few rows only)
$ws.Range("F7").Value  = "Reference(Patient)"
$ws.Range("F9").Value  = "Reference(Practitioner | Device)"
$ws.Range("F16").Value = "Reference(Practitioner) - Who is committing?"

# Column B - HL7 v2 Mapping as per Argonaut/www.hl7.org
$ws.Range("B3").Value  = "IAM-7"
$ws.Range("B5").Value  = "AL1-2 / IAM-9"
$ws.Range("B6").Value  = "AL1-4/ IAR-2"
$ws.Range("B7").Value  = "PID-3"
$ws.Range("B8").Value  = "ALI-6 / IAM-11"
$ws.Range("B9").Value  = "IAM-14 / IAM-16"
$ws.Range("B11").Value = "AL1-5 / IAM-8"

# Column C - HL7 Field (Numeric)
$ws.Range("C3").Value  = "IAM.7"
$ws.Range("C5").Value  = "AL1.2 || IAM.9"
$ws.Range("C6").Value  = "AL1.4|| IAR.2"
$ws.Range("C7").Value  = "PID.3"
$ws.Range("C8").Value  = "ALI.6 || IAM.11"
$ws.Range("C9").Value  = "IAM.14 || IAM.16"
$ws.Range("C11").Value = "AL1.5 || IAM.8"

# Column D - HL7 Field Name
$ws.Range("D3").Value  = "IAM.AllergyUniqueIdentifier"
$ws.Range("D5").Value  = "AL1.AllergenTypeCode || IAM.SensitivitytoCausativeAgentC"
$ws.Range("D6").Value  = "AL1.AllergySeverityCode || IAR.AllergySeverityCode"
$ws.Range("D7").Value  = "PID.PatientIdentifierList"
$ws.Range("D8").Value  = "AL1.IdentificationDate || IAM.OnsetDate"
$ws.Range("D9").Value  = "IAM.ReportedBy || IAM.AlertDeviceCode"
$ws.Range("D11").Value = "AL1.AllergyReactionCode || IAM.ActionReason"

# Match the author's final cursor position (cell D7 selected).
$ws.Range("D7").Select()
